$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mis-decoded accented characters in the Regional Economic Communities note
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# Update data values in row 97
$ws.Range("C97").Value = 349509.3
$ws.Range("D97").Value = 14999.3
$ws.Range("E97").Value = 35106.699999999997
$ws.Range("F97").Value = 140551.79999999999
$ws.Range("G97").Value = 87872.3
$ws.Range("H97").Value = 5411.9
$ws.Range("I97").Value = 2890.9
$ws.Range("K97").Value = 2865
$ws.Range("L97").Value = 137867.29999999999
$ws.Range("M97").Value = 59695.9

# Update data values in row 98
$ws.Range("C98").Value = 1119955.2
$ws.Range("D98").Value = 64598.6
$ws.Range("E98").Value = 218295.9
$ws.Range("F98").Value = 264579.8
$ws.Range("G98").Value = 168529.9
$ws.Range("H98").Value = 39406.1
$ws.Range("I98").Value = 20501.7
$ws.Range("K98").Value = 19302.599999999999
$ws.Range("L98").Value = 244638.2
$ws.Range("M98").Value = 324618.90000000002
